$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four oldest year rows (2002-2005); this shifts 2016-2020 up
# into rows 2-6, matching the new layout. (Deleting one row at a time —
# multi-row range deletes misbehave in this COM host.)
$ws.Range("A2:B2").Delete()
$ws.Range("A2:B2").Delete()
$ws.Range("A2:B2").Delete()
$ws.Range("A2:B2").Delete()

# Add the new 2021 row at the end (now row 7), matching the formatting
# used by the row above it (bold, centered, bordered label style).
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = "2021年"
$ws.Range("B7").Value = 3
